$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 21:45"

# Country label reorders (Column A) - ranking shuffled by the new day's data
$ws.Range("A92").Value = "Mauritania"
$ws.Range("A93").Value = "Hungria"
$ws.Range("A126").Value = "Cabo Verde"
$ws.Range("A127").Value = "Benin"
$ws.Range("A128").Value = "Malaui"
$ws.Range("A129").Value = "Jordania"
$ws.Range("A130").Value = "Yemen"
$ws.Range("A131").Value = "Letonia"
$ws.Range("A143").Value = "Suazilandia"
$ws.Range("A144").Value = "Liberia"
$ws.Range("A171").Value = "Namibia"
$ws.Range("A172").Value = "Guadalupe"
$ws.Range("A173").Value = "Gibraltar"
$ws.Range("A174").Value = "Burundi"
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# Updated statistics (columns B-H), one cell at a time
$ws.Range("B4").Value = 2621662
$ws.Range("C4").Value = 25125
$ws.Range("D4").Value = 1083709
$ws.Range("E4").Value = 1409638
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 163
$ws.Range("H4").Value = 128315

$ws.Range("B5").Value = 1323069
$ws.Range("C5").Value = 7128
$ws.Range("D5").Value = 715905
$ws.Range("E5").Value = 549990
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 71
$ws.Range("H5").Value = 57174

$ws.Range("B7").Value = 549197
$ws.Range("C7").Value = 19620
$ws.Range("D7").Value = 321774
$ws.Range("E7").Value = 210936
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 384
$ws.Range("H7").Value = 16487

$ws.Range("B16").Value = 197239
$ws.Range("C16").Value = 1356
$ws.Range("D16").Value = 170595
$ws.Range("E16").Value = 21547
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 5097

$ws.Range("B68").Value = 12052
$ws.Range("C68").Value = 175
$ws.Range("D68").Value = 8740
$ws.Range("E68").Value = 3091
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 221

$ws.Range("B92").Value = 4149
$ws.Range("C92").Value = 124
$ws.Range("D92").Value = 1419
$ws.Range("E92").Value = 2604
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 5
$ws.Range("H92").Value = 126

$ws.Range("B93").Value = 4142
$ws.Range("C93").Value = 4
$ws.Range("D93").Value = 2685
$ws.Range("E93").Value = 876
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = 581

$ws.Range("B121").Value = 1545
$ws.Range("C121").Value = 14
$ws.Range("D121").Value = 1289
$ws.Range("E121").Value = 234
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 22

$ws.Range("B126").Value = 1155
$ws.Range("C126").Value = 64
$ws.Range("D126").Value = 570
$ws.Range("E126").Value = 573
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 12

$ws.Range("B127").Value = 1149
$ws.Range("C127").Value = 25
$ws.Range("D127").Value = 306
$ws.Range("E127").Value = 827
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 2
$ws.Range("H127").Value = 16

$ws.Range("B128").Value = 1146
$ws.Range("C128").Value = 108
$ws.Range("D128").Value = 260
$ws.Range("E128").Value = 873
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 13

$ws.Range("B129").Value = 1121
$ws.Range("C129").Value = 10
$ws.Range("D129").Value = 860
$ws.Range("E129").Value = 252
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 9

$ws.Range("B130").Value = 1118
$ws.Range("C130").Value = 15
$ws.Range("D130").Value = 430
$ws.Range("E130").Value = 386
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 6
$ws.Range("H130").Value = 302

$ws.Range("B131").Value = 1116
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 932
$ws.Range("E131").Value = 154
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 30

$ws.Range("B135").Value = 959
$ws.Range("C135").Value = 18
$ws.Range("D135").Value = 830
$ws.Range("E135").Value = 76
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 53

$ws.Range("B138").Value = 900
$ws.Range("C138").Value = 22
$ws.Range("D138").Value = 443
$ws.Range("E138").Value = 455
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 2

$ws.Range("B143").Value = 781
$ws.Range("C143").Value = 36
$ws.Range("D143").Value = 372
$ws.Range("E143").Value = 398
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 3
$ws.Range("H143").Value = 11

$ws.Range("B144").Value = 768
$ws.Range("C144").Value = 39
$ws.Range("D144").Value = 298
$ws.Range("E144").Value = 436
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 34

$ws.Range("B153").Value = 521
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 472
$ws.Range("E153").Value = 47
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 2

$ws.Range("B171").Value = 183
$ws.Range("C171").Value = 47
$ws.Range("D171").Value = 24
$ws.Range("E171").Value = 159
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

$ws.Range("B172").Value = 182
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 157
$ws.Range("E172").Value = 11
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 14

$ws.Range("B173").Value = 177
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 176
$ws.Range("E173").Value = 1
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

$ws.Range("B174").Value = 170
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 115
$ws.Range("E174").Value = 54
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 1

$ws.Range("B196").Value = 27
$ws.Range("C196").Value = 3
$ws.Range("D196").Value = 4
$ws.Range("E196").Value = 23
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

Write-Output "done"